$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 366, shifting existing rows 366..445 down to 367..446
$ws.Rows(366).Insert()

# Populate the new row 366 with the new record
$ws.Cells.Item(366, 1).Value = 9
$ws.Cells.Item(366, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(366, 3).Value = "Metropolitana"
$ws.Cells.Item(366, 4).Value = 45173
$ws.Cells.Item(366, 5).Value = 13
$ws.Cells.Item(366, 6).Value = 100112021
$ws.Cells.Item(366, 7).Value = "Ají"
$ws.Cells.Item(366, 8).Value = "Inferno"
$ws.Cells.Item(366, 9).Value = "Primera"
$ws.Cells.Item(366, 10).Value = 70
$ws.Cells.Item(366, 11).Value = 14000
$ws.Cells.Item(366, 12).Value = 15000
$ws.Cells.Item(366, 13).Value = 14500
$ws.Cells.Item(366, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(366, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(366, 16).Value = 1450
$ws.Cells.Item(366, 17).Value = 10
$ws.Cells.Item(366, 18).Value = "Hortaliza"
